$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-8 from 2023-11-13 (45243)
# to 2023-11-14 (45244).
foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45244
}
